# Apply the "Deploying to gh-pages ... LinuxForHealth" update to the
# StructureDefinition-communication-mode workbook.
#
# Sheet "Metadata": rebrand IBM/Alvearie references to LinuxForHealth,
#   bump the version and republish date.
# Sheet "Elements": rebrand the bound ValueSet URL, and move the
#   ele-1/ext-1 constraint note off the "Extension" row and onto the
#   "Extension.extension" row where it really belongs.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/communication-mode"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Y7").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/eng-communication-mode"

$constraint = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
$elements.Range("AI2").Value = ""
$elements.Range("AI4").Value = $constraint
